$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Title cell
$ws.Range("B1").Value2 = "Research Assignment Planner"

# 2. Activity row 5: "Catchup with all chapters" -> "Conduct Indivisual Research"; D5/F5: 3 -> 1
$ws.Range("B5").Value2 = "Conduct Indivisual Research"
$ws.Range("D5").Value2 = 1
$ws.Range("F5").Value2 = 1

# 3. Activity row 6: "Project Planning" -> "Group Discussion"; C6/E6: 1 -> 2
$ws.Range("B6").Value2 = "Group Discussion"
$ws.Range("C6").Value2 = 2
$ws.Range("E6").Value2 = 2

# 4. Activity row 7: "Formulating Mysteries" -> "Formulating Slides"; C7/E7: 4 -> 2
$ws.Range("B7").Value2 = "Formulating Slides"
$ws.Range("C7").Value2 = 2
$ws.Range("E7").Value2 = 2

# 5. Activity row 8: "Functional Specification" -> "Film & Edit Video Explanation"; C8/E8: 4 -> 3
$ws.Range("B8").Value2 = "Film & Edit Video Explanation"
$ws.Range("C8").Value2 = 3
$ws.Range("E8").Value2 = 3

# 6. Rows 9-15 no longer hold activities - clear their content/formatting entirely
$ws.Range("B9:AA15").ClearContents()
$ws.Range("B9:AA15").ClearFormats()
# rows 14 & 15 had custom heights (26 / 33) - restore the sheet's default of 30
$ws.Rows("14:15").RowHeight = 30

# 7. The trailing blank rows 108-114 are removed entirely, shortening the sheet to row 107
$ws.Rows("108:114").Delete()

# 8. The conditional formatting that covered the (now empty) activity rows 5-15
#    only covers rows 5-8 now
$fcs = $ws.Range("H5:BO15").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("H5:BO8"))
}

# 9. Selection moved to F11
$ws.Range("F11").Select()
